$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 7064.5
$ws.Cells.Item(17, 10).Value = 8119.294
$ws.Cells.Item(17, 12).Value = 24357.882
$ws.Cells.Item(17, 14).Value = -24693.882
$ws.Cells.Item(33, 8).Value = 122.63636
$ws.Cells.Item(33, 9).Value = 127.888885
$ws.Cells.Item(33, 10).Value = 99
$ws.Cells.Item(33, 11).Value = 127.888885
$ws.Cells.Item(33, 12).Value = 99
$ws.Cells.Item(33, 13).Value = 101.111115
$ws.Cells.Item(33, 14).Value = -557
$ws.Cells.Item(62, 8).Value = 3905.9443
$ws.Cells.Item(62, 9).Value = 3050.625
$ws.Cells.Item(62, 10).Value = 4590.2
$ws.Cells.Item(62, 11).Value = 3050.625
$ws.Cells.Item(62, 12).Value = 4590.2
$ws.Cells.Item(62, 13).Value = -2426.625
$ws.Cells.Item(62, 14).Value = -5838.2
$ws.Cells.Item(64, 8).Value = 4312.375
$ws.Cells.Item(64, 9).Value = 3749.5
$ws.Cells.Item(64, 11).Value = 3749.5
$ws.Cells.Item(64, 13).Value = -3501.5
$ws.Cells.Item(65, 8).Value = 3905.9443
$ws.Cells.Item(65, 9).Value = 3050.625
$ws.Cells.Item(65, 10).Value = 4590.2
$ws.Cells.Item(65, 11).Value = 15253.125
$ws.Cells.Item(65, 12).Value = 22951
$ws.Cells.Item(65, 13).Value = -12133.125
$ws.Cells.Item(65, 14).Value = -29191
$ws.Cells.Item(67, 8).Value = 4312.375
$ws.Cells.Item(67, 9).Value = 3749.5
$ws.Cells.Item(67, 11).Value = 3749.5
$ws.Cells.Item(67, 13).Value = -2891.5
$ws.Cells.Item(86, 8).Value = 7896.294
$ws.Cells.Item(86, 10).Value = 18332
$ws.Cells.Item(86, 12).Value = 18332
$ws.Cells.Item(86, 14).Value = -20578
$ws.Cells.Item(89, 8).Value = 7896.294
$ws.Cells.Item(89, 10).Value = 18332
$ws.Cells.Item(89, 12).Value = 91660
$ws.Cells.Item(89, 14).Value = -102892
$ws.Cells.Item(129, 8).Value = 1086.4916
$ws.Cells.Item(129, 10).Value = 1186.8431
$ws.Cells.Item(129, 12).Value = 3560.5293
$ws.Cells.Item(129, 14).Value = -13560.5293
$ws.Cells.Item(132, 8).Value = 2612.389
$ws.Cells.Item(132, 9).Value = 2710.1765
$ws.Cells.Item(132, 10).Value = 950
$ws.Cells.Item(132, 11).Value = 8130.529500000001
$ws.Cells.Item(132, 12).Value = 2850
$ws.Cells.Item(132, 13).Value = -5600.529500000001
$ws.Cells.Item(132, 14).Value = -7910
$ws.Cells.Item(138, 8).Value = 2365.1206
$ws.Cells.Item(138, 9).Value = 2195.8
$ws.Cells.Item(138, 10).Value = 2424.186
$ws.Cells.Item(138, 11).Value = 6587.400000000001
$ws.Cells.Item(138, 12).Value = 7272.558000000001
$ws.Cells.Item(138, 13).Value = -1447.400000000001
$ws.Cells.Item(138, 14).Value = -17552.558

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5725.2964
$ws.Cells.Item(32, 9).Value = 4429.2085
$ws.Cells.Item(32, 11).Value = 4429.2085
$ws.Cells.Item(32, 13).Value = -4142.2085
$ws.Cells.Item(61, 8).Value = 1214.5834
$ws.Cells.Item(61, 9).Value = 1049.3135
$ws.Cells.Item(61, 10).Value = 3429.2
$ws.Cells.Item(61, 11).Value = 1049.3135
$ws.Cells.Item(61, 12).Value = 3429.2
$ws.Cells.Item(61, 13).Value = -837.3135
$ws.Cells.Item(61, 14).Value = -3853.2
$ws.Cells.Item(110, 8).Value = 478.0909
$ws.Cells.Item(110, 9).Value = 469.8889
$ws.Cells.Item(110, 10).Value = 515
$ws.Cells.Item(110, 11).Value = 469.8889
$ws.Cells.Item(110, 12).Value = 515
$ws.Cells.Item(110, 13).Value = 1575.1111
$ws.Cells.Item(110, 14).Value = -4605
$ws.Cells.Item(132, 8).Value = 15802.194
$ws.Cells.Item(132, 9).Value = 1809.0646
$ws.Cells.Item(132, 10).Value = 102559.6
$ws.Cells.Item(132, 11).Value = 5427.1938
$ws.Cells.Item(132, 12).Value = 307678.8
$ws.Cells.Item(132, 13).Value = -2897.1938
$ws.Cells.Item(132, 14).Value = -312738.8
$ws.Cells.Item(136, 8).Value = 1214.5834
$ws.Cells.Item(136, 9).Value = 1049.3135
$ws.Cells.Item(136, 10).Value = 3429.2
$ws.Cells.Item(136, 11).Value = 3147.9405
$ws.Cells.Item(136, 12).Value = 10287.6
$ws.Cells.Item(136, 13).Value = -597.9404999999997
$ws.Cells.Item(136, 14).Value = -15387.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3505.0588
$ws.Cells.Item(134, 9).Value = 3505.0588
$ws.Cells.Item(134, 11).Value = 10515.1764
$ws.Cells.Item(134, 13).Value = -7980.1764

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 13668.134
$ws.Cells.Item(31, 9).Value = 16468.096
$ws.Cells.Item(31, 11).Value = 16468.096
$ws.Cells.Item(31, 13).Value = -16173.096
$ws.Cells.Item(34, 8).Value = 13668.134
$ws.Cells.Item(34, 9).Value = 16468.096
$ws.Cells.Item(34, 11).Value = 16468.096
$ws.Cells.Item(34, 13).Value = -16266.096
$ws.Cells.Item(58, 8).Value = 12787.619
$ws.Cells.Item(58, 9).Value = 866.25
$ws.Cells.Item(58, 10).Value = 36630.355
$ws.Cells.Item(58, 11).Value = 866.25
$ws.Cells.Item(58, 12).Value = 36630.355
$ws.Cells.Item(58, 13).Value = -663.25
$ws.Cells.Item(58, 14).Value = -37036.355
$ws.Cells.Item(100, 8).Value = 170390
$ws.Cells.Item(100, 10).Value = 170390
$ws.Cells.Item(100, 12).Value = 170390
$ws.Cells.Item(100, 14).Value = -172554
$ws.Cells.Item(105, 8).Value = 62502504
$ws.Cells.Item(105, 9).Value = 125000000
$ws.Cells.Item(105, 11).Value = 125000000
$ws.Cells.Item(105, 13).Value = -124998253
$ws.Cells.Item(127, 8).Value = 25325.834
$ws.Cells.Item(127, 9).Value = 15000
$ws.Cells.Item(127, 10).Value = 39782
$ws.Cells.Item(127, 11).Value = 15000
$ws.Cells.Item(127, 12).Value = 39782
$ws.Cells.Item(127, 13).Value = -10040
$ws.Cells.Item(127, 14).Value = -49702
$ws.Cells.Item(132, 8).Value = 14132.419
$ws.Cells.Item(132, 9).Value = 18575.033
$ws.Cells.Item(132, 10).Value = 3880.2307
$ws.Cells.Item(132, 11).Value = 55725.099
$ws.Cells.Item(132, 12).Value = 11640.6921
$ws.Cells.Item(132, 13).Value = -53195.099
$ws.Cells.Item(132, 14).Value = -16700.6921
$ws.Cells.Item(134, 8).Value = 794.0526
$ws.Cells.Item(134, 9).Value = 722.7646999999999
$ws.Cells.Item(134, 10).Value = 1400
$ws.Cells.Item(134, 11).Value = 2168.2941
$ws.Cells.Item(134, 12).Value = 4200
$ws.Cells.Item(134, 13).Value = 366.7058999999999
$ws.Cells.Item(134, 14).Value = -9270
$ws.Cells.Item(136, 8).Value = 12787.619
$ws.Cells.Item(136, 9).Value = 866.25
$ws.Cells.Item(136, 10).Value = 36630.355
$ws.Cells.Item(136, 11).Value = 2598.75
$ws.Cells.Item(136, 12).Value = 109891.065
$ws.Cells.Item(136, 13).Value = -48.75
$ws.Cells.Item(136, 14).Value = -114991.065

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 455.41666
$ws.Cells.Item(122, 9).Value = 242.47058
$ws.Cells.Item(122, 11).Value = 2182.23522
$ws.Cells.Item(122, 13).Value = 267.76478

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 25691.709
$ws.Cells.Item(132, 9).Value = 5155.15
$ws.Cells.Item(132, 10).Value = 128374.5
$ws.Cells.Item(132, 11).Value = 15465.45
$ws.Cells.Item(132, 12).Value = 385123.5
$ws.Cells.Item(132, 13).Value = -12935.45
$ws.Cells.Item(132, 14).Value = -390183.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4953.12
$ws.Cells.Item(7, 9).Value = 3423.4285
$ws.Cells.Item(7, 10).Value = 6900
$ws.Cells.Item(7, 11).Value = 3423.4285
$ws.Cells.Item(7, 12).Value = 6900
$ws.Cells.Item(7, 13).Value = -3311.4285
$ws.Cells.Item(7, 14).Value = -7124
$ws.Cells.Item(64, 8).Value = 22481
$ws.Cells.Item(64, 10).Value = 22481
$ws.Cells.Item(64, 12).Value = 22481
$ws.Cells.Item(64, 14).Value = -22931
$ws.Cells.Item(67, 8).Value = 22481
$ws.Cells.Item(67, 10).Value = 22481
$ws.Cells.Item(67, 12).Value = 22481
$ws.Cells.Item(67, 14).Value = -24041
$ws.Cells.Item(68, 8).Value = 2136.5833
$ws.Cells.Item(68, 9).Value = 1454.875
$ws.Cells.Item(68, 10).Value = 3500
$ws.Cells.Item(68, 11).Value = 1454.875
$ws.Cells.Item(68, 12).Value = 3500
$ws.Cells.Item(68, 13).Value = -705.875
$ws.Cells.Item(68, 14).Value = -4998
$ws.Cells.Item(71, 8).Value = 2136.5833
$ws.Cells.Item(71, 9).Value = 1454.875
$ws.Cells.Item(71, 10).Value = 3500
$ws.Cells.Item(71, 11).Value = 7274.375
$ws.Cells.Item(71, 12).Value = 17500
$ws.Cells.Item(71, 13).Value = -3530.375
$ws.Cells.Item(71, 14).Value = -24988
$ws.Cells.Item(82, 8).Value = 3263.9092
$ws.Cells.Item(82, 10).Value = 3000.5
$ws.Cells.Item(82, 12).Value = 3000.5
$ws.Cells.Item(82, 14).Value = -3722.5
$ws.Cells.Item(85, 8).Value = 3263.9092
$ws.Cells.Item(85, 10).Value = 3000.5
$ws.Cells.Item(85, 12).Value = 3000.5
$ws.Cells.Item(85, 14).Value = -5496.5
$ws.Cells.Item(87, 8).Value = 12800
$ws.Cells.Item(87, 9).Value = 12800
$ws.Cells.Item(87, 11).Value = 12800
$ws.Cells.Item(87, 13).Value = -11677
$ws.Cells.Item(88, 8).Value = 38000
$ws.Cells.Item(88, 10).Value = 38000
$ws.Cells.Item(88, 12).Value = 38000
$ws.Cells.Item(88, 14).Value = -38856
$ws.Cells.Item(90, 8).Value = 12800
$ws.Cells.Item(90, 9).Value = 12800
$ws.Cells.Item(90, 11).Value = 38400
$ws.Cells.Item(90, 13).Value = -32784
$ws.Cells.Item(91, 8).Value = 38000
$ws.Cells.Item(91, 10).Value = 38000
$ws.Cells.Item(91, 12).Value = 38000
$ws.Cells.Item(91, 14).Value = -40964
$ws.Cells.Item(126, 8).Value = 4953.12
$ws.Cells.Item(126, 9).Value = 3423.4285
$ws.Cells.Item(126, 10).Value = 6900
$ws.Cells.Item(126, 11).Value = 10270.2855
$ws.Cells.Item(126, 12).Value = 20700
$ws.Cells.Item(126, 13).Value = -7800.2855
$ws.Cells.Item(126, 14).Value = -25640
$ws.Cells.Item(132, 8).Value = 2038.0769
$ws.Cells.Item(132, 9).Value = 1405.875
$ws.Cells.Item(132, 11).Value = 4217.625
$ws.Cells.Item(132, 13).Value = -1687.625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(63, 8).Value = 39999
$ws.Cells.Item(63, 10).Value = 39999
$ws.Cells.Item(63, 12).Value = 39999
$ws.Cells.Item(63, 14).Value = -41247
$ws.Cells.Item(66, 8).Value = 39999
$ws.Cells.Item(66, 10).Value = 39999
$ws.Cells.Item(66, 12).Value = 119997
$ws.Cells.Item(66, 14).Value = -126237
$ws.Cells.Item(69, 8).Value = 11800.4
$ws.Cells.Item(69, 10).Value = 13750.25
$ws.Cells.Item(69, 12).Value = 13750.25
$ws.Cells.Item(69, 14).Value = -15248.25
$ws.Cells.Item(72, 8).Value = 11800.4
$ws.Cells.Item(72, 10).Value = 13750.25
$ws.Cells.Item(72, 12).Value = 41250.75
$ws.Cells.Item(72, 14).Value = -48738.75
$ws.Cells.Item(122, 8).Value = 1781.8462
$ws.Cells.Item(122, 9).Value = 1781.8462
$ws.Cells.Item(122, 11).Value = 5345.5386
$ws.Cells.Item(122, 13).Value = -2895.5386
$ws.Cells.Item(126, 8).Value = 1076.5
$ws.Cells.Item(126, 9).Value = 1019
$ws.Cells.Item(126, 10).Value = 1249
$ws.Cells.Item(126, 11).Value = 3057
$ws.Cells.Item(126, 12).Value = 3747
$ws.Cells.Item(126, 13).Value = -587
$ws.Cells.Item(126, 14).Value = -8687
$ws.Cells.Item(132, 8).Value = 929.63043
$ws.Cells.Item(132, 9).Value = 667.75757
$ws.Cells.Item(132, 10).Value = 1594.3846
$ws.Cells.Item(132, 11).Value = 2003.27271
$ws.Cells.Item(132, 12).Value = 4783.1538
$ws.Cells.Item(132, 13).Value = 526.72729
$ws.Cells.Item(132, 14).Value = -9843.1538
$ws.Cells.Item(136, 8).Value = 35715892
$ws.Cells.Item(136, 9).Value = 43479730
$ws.Cells.Item(136, 11).Value = 130439190
$ws.Cells.Item(136, 13).Value = -130436640
